$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.411.63"
$ws.Range("E2").Value = "  +0.05%  "
$ws.Range("D3").Value = "1.848.53"
$ws.Range("E3").Value = "  +0.21%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9995"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "240.52"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.68%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6282"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.51%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9999"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07706"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +2.28%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2921"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -0.32%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "25.03"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +1.97%  "
$ws.Range("E11").Value = "  +0.44%  "
$ws.Range("D12").Value = "1.857.21"
$ws.Range("E12").Value = "  +0.81%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.00001086"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +4.09%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6827"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +0.40%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "83.56"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +0.27%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.191"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +0.32%  "
$ws.Range("D18").Value = "29.430.83"
$ws.Range("E18").Value = "  +0.01%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "229.02"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +0.10%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.39"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -0.14%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.9998"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -0.06%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.461"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +0.04%  "
$ws.Range("E23").Value = "  -0.03%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "157.60"
$ws.Range("D24").ClearFormats()
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1379"
$ws.Range("D25").ClearFormats()
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.425"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +0.86%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "17.70"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +0.64%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.352"
$ws.Range("D28").ClearFormats()
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.461"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +0.24%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05639"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +0.17%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.125"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +0.49%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.047"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +0.55%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.843"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -0.21%  "
$ws.Range("E34").Value = "  +0.73%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7028"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -1.37%  "
$ws.Range("E36").Value = "  +0.09%  "
$ws.Range("D37").Value = "1.226.34"
$ws.Range("E37").Value = "  -1.64%  "
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01790"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -1.11%  "
$ws.Range("B39").Value = "MXToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.756"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -0.38%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.447"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +1.20%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9075"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +0.58%  "
$ws.Range("B42").Value = "RocketPoolETH"
$ws.Range("C42").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D42").Value = "2.043.78"
$ws.Range("E42").Value = "  +2.33%  "
$ws.Range("B43").Value = "PaxDollar"
$ws.Range("C43").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.000"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +0.02%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "101.84"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +0.14%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "66.09"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +0.43%  "
$ws.Range("E46").Value = "  +2.26%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.188"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +1.20%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4020"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +0.62%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.1158"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +3.21%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.008"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +0.95%  "
$ws.Range("E51").Value = "  +0.39%  "
